$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.643.27"
$ws.Range("E2").Value = "  +3.44%  "
$ws.Range("D3").Value = "1.700.62"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.97"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4040"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.543"
$ws.Range("E9").Value = "  +8.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.87"
$ws.Range("E10").Value = "  +13.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.000"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08824"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.306"
$ws.Range("E13").Value = "  +8.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.44"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001335"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.646"
$ws.Range("E16").Value = "  +6.35%  "
$ws.Range("D17").Value = "1.710.02"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.82"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07068"
$ws.Range("E19").Value = "  +4.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.73"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.948"
$ws.Range("E21").Value = "  +4.77%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.18"
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("D24").Value = "24.630.34"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.987"
$ws.Range("E25").Value = "  +9.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.328"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.45"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.20"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.231"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.03"
$ws.Range("E30").Value = "  +3.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.780"
$ws.Range("E31").Value = "  +18.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.115"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").Value = "1.887.99"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.455"
$ws.Range("E34").Value = "  +13.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08582"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.20"
$ws.Range("E36").Value = "  +7.86%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2767"
$ws.Range("E37").Value = "  +4.98%  "
$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.958"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.81"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02785"
$ws.Range("E40").Value = "  +10.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09055"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.472"
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7789"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7305"
$ws.Range("E44").Value = "  +4.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.77"
$ws.Range("E45").Value = "  +6.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.519"
$ws.Range("E46").Value = "  +6.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.209"
$ws.Range("E47").Value = "  +3.45%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.08"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.336"
$ws.Range("E50").Value = "  +16.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08044"
$ws.Range("E51").Value = "  +3.79%  "
